$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at 62 (pushes the old row 62, which only has D62, down to row 63)
$ws.Rows.Item(62).Insert()

# Fill in the new row 62 data (leading apostrophe matches the quote-prefix
# formatting used by the other "extension" rows in column A, e.g. A60/A61)
$ws.Cells.Item(62, 1).Value = "'     -"
$ws.Cells.Item(62, 2).Value = "US Core Tribal Affiliation Extension"
$ws.Cells.Item(62, 3).Value = "6.0.0"

# Update the sheet view to match the new scroll/selection state
$excel.ActiveWindow.ScrollRow = 40
$ws.Range("A65").Select()
